$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 15 and 16 as reserved/bought ("Y") in column E
$ws.Range("E15").Value = "Y"
$ws.Range("E16").Value = "Y"

# New row 17: Puzzle gift idea
$ws.Range("A17").Value = "Puzzle (orice 2-5 ani)"
$ws.Range("B17").Value = "https://assets.dryicons.com/uploads/icon/svg/13100/1a36a162-9ce3-4d75-a031-b61625848a32.svg"

# New row 18: Zalando gift card
$ws.Range("A18").Value = "Card cadou Zalando"
$ws.Range("C18").Value = "https://zalando.ch"
$ws.Range("B18").Value = "https://img01.ztat.net/article/spp-media-p1/bfd222ac1b2541db8553ad0bfbde422d/54a998916cfc4620b4bac8119c09a041.jpg?imwidth=300&filter=packshot"

# Update the active selection to match the author's final cursor position
$ws.Range("E18").Select()
